$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("DAVILA CORDOVA MARIBEL", 45),
    @("FARFAN MONTOYA ROSSANA ISABEL", 44),
    @("ANGIE BELÉN RODRÍGUEZ ZAVALA", 42),
    @("SANCHEZ ULLOA CESAR AUGUSTO", 42),
    @("FIORELA KEILY GUTIERREZ CRUZ", 41),
    @("RAMOS RAMOS HANDY JAIR", 41),
    @("GONZALES VALLE SEBASTIAN", 40),
    @("GARCIA GUTIERREZ LUIS ARTURO", 39),
    @("CASTILLO QUEZADA DIEGO ALONSO", 38),
    @("OLIVA ALVA GOSSELYN NASSIRA", 37),
    @("ARANEDA LOPEZ MARCO VIERI", 37),
    @("DELGADO DELGADO RONI", 37),
    @("CARBAJAL RAMOS JESUS MARINA", 37),
    @("DANY DARWIN VILLACORTA SAAVEDRA", 35),
    @("VERDE LIZARRAGA DEYSI EUFEMIA", 35),
    @("JOSSY IVANA SUÁREZ ZAVALETA", 35),
    @("BAZAN TEJADA JOSE VICENTE", 34),
    @("GUERRA CALDERON ESTHEFANY NICOLLE", 33),
    @("JAVE CHAVEZ ANGHELO MARTIN", 32),
    @("CARDENAS CAMPOJO MARY PAULA", 28),
    @("ALVITES CAMPOS SERGIO MARTIN", 13),
    @("ROCHA SIPIRAN JHORDAN ENRIQUE", 1)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
